# Update column F (dSF) values per the re-pulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -3
    3  = -4
    4  = -3
    5  = -1
    6  = -5
    7  = 5
    8  = 7
    9  = 6
    10 = -1
    12 = 3
    13 = -5
    14 = 8
    16 = -3
    17 = -2
    18 = -2
    19 = 1
    20 = 2
    21 = -2
    22 = 1
    23 = 2
    24 = -4
    25 = -5
    26 = 3
    27 = -4
    28 = -6
    29 = -3
    30 = 6
    31 = 2
    32 = 1
    33 = 3
    34 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
